{"js": "// Update the date line and the twenty-five \"XXX\u00f7Y=\" division problems to\n// the new values from the next day's worksheet.\nconst replacements = [\n  [\"2024-09-28 Saturday\", \"2024-09-29 Sunday\"],\n  [\"866\u00f74=\", \"573\u00f77=\"],\n  [\"816\u00f77=\", \"607\u00f72=\"],\n  [\"687\u00f73=\", \"337\u00f79=\"],\n  [\"936\u00f77=\", \"405\u00f72=\"],\n  [\"754\u00f72=\", \"295\u00f73=\"],\n  [\"285\u00f77=\", \"187\u00f78=\"],\n  [\"159\u00f74=\", \"854\u00f78=\"],\n  [\"375\u00f76=\", \"133\u00f76=\"],\n  [\"177\u00f79=\", \"186\u00f77=\"],\n  [\"440\u00f77=\", \"145\u00f75=\"],\n  [\"934\u00f76=\", \"210\u00f73=\"],\n  [\"486\u00f72=\", \"251\u00f76=\"],\n  [\"588\u00f78=\", \"706\u00f79=\"],\n  [\"870\u00f77=\", \"532\u00f72=\"],\n  [\"108\u00f72=\", \"998\u00f77=\"],\n  [\"680\u00f74=\", \"838\u00f77=\"],\n  [\"717\u00f74=\", \"841\u00f77=\"],\n  [\"160\u00f74=\", \"792\u00f75=\"],\n  [\"904\u00f76=\", \"230\u00f76=\"],\n  [\"826\u00f79=\", \"913\u00f78=\"],\n  [\"811\u00f72=\", \"108\u00f77=\"],\n  [\"259\u00f74=\", \"211\u00f74=\"],\n  [\"991\u00f73=\", \"689\u00f76=\"],\n  [\"374\u00f74=\", \"409\u00f75=\"],\n  [\"255\u00f74=\", \"955\u00f77=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and the twenty-five \"XXX\u00f7Y=\" division problems to\n# the new values from the next day's worksheet.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-09-28 Saturday\", \"2024-09-29 Sunday\"),\n    @(\"866\u00f74=\", \"573\u00f77=\"),\n    @(\"816\u00f77=\", \"607\u00f72=\"),\n    @(\"687\u00f73=\", \"337\u00f79=\"),\n    @(\"936\u00f77=\", \"405\u00f72=\"),\n    @(\"754\u00f72=\", \"295\u00f73=\"),\n    @(\"285\u00f77=\", \"187\u00f78=\"),\n    @(\"159\u00f74=\", \"854\u00f78=\"),\n    @(\"375\u00f76=\", \"133\u00f76=\"),\n    @(\"177\u00f79=\", \"186\u00f77=\"),\n    @(\"440\u00f77=\", \"145\u00f75=\"),\n    @(\"934\u00f76=\", \"210\u00f73=\"),\n    @(\"486\u00f72=\", \"251\u00f76=\"),\n    @(\"588\u00f78=\", \"706\u00f79=\"),\n    @(\"870\u00f77=\", \"532\u00f72=\"),\n    @(\"108\u00f72=\", \"998\u00f77=\"),\n    @(\"680\u00f74=\", \"838\u00f77=\"),\n    @(\"717\u00f74=\", \"841\u00f77=\"),\n    @(\"160\u00f74=\", \"792\u00f75=\"),\n    @(\"904\u00f76=\", \"230\u00f76=\"),\n    @(\"826\u00f79=\", \"913\u00f78=\"),\n    @(\"811\u00f72=\", \"108\u00f77=\"),\n    @(\"259\u00f74=\", \"211\u00f74=\"),\n    @(\"991\u00f73=\", \"689\u00f76=\"),\n    @(\"374\u00f74=\", \"409\u00f75=\"),\n    @(\"255\u00f74=\", \"955\u00f77=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
